# "Iteration Burndown" sheet: the iteration's scope changed -- the team's
# starting "Hours Left" dropped from 13 to 11, and every day after day 0
# now reads 0 (burndown reset for the remaining tracked days). Column C
# ("Ideal") is a formula column (=B2-(B2/14)*A#) and recalculates on its
# own once B2 changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11
$ws.Range("B3:B15").Value = 0

# Leave the selection where the edit happened.
$ws.Range("B3").Select() | Out-Null
